$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A105").Value = "신멸인(燼滅刃)"
$ws.Range("B105").Value = "신멸인의 혼(燼滅刃の魂)"
$ws.Range("D105").Value = "예리도 레벨+2(斬れ味レベル+2) + 탄도강화(弾導強化) + 숫돌 사용 고속화(砥石使用高速化)"

$ws.Range("A106").Value = "진・개열(真・鎧裂)"
$ws.Range("B106").Value = "진・개열의 혼(真・鎧裂の魂)"
$ws.Range("D106").Value = "가드 성능+2(ガード性能+2) + 납도술(納刀術) + 인린연마(刃鱗磨き)"

$ws.Range("A107").Value = "인린(刃鱗)"
$ws.Range("B107").Value = "인린연마(刃鱗磨き)"
$ws.Range("D107").Value = "발도중에 회피행동에 의해 효과발동. 검사는 예리도 5회복(셀레기오스 무기는 7회복). 보우건은 탄 1발 장전. 50%의 확률로 장전수+1. 활은 접격병의 데미지가 1.5배(셀레기오스 무기는 1.65배)"

$ws.Range("A108").Value = "수면(睡眠)"
$ws.Range("B108").Value = "수면무효(睡眠無効)/수면배가(睡眠倍加)"
$ws.Range("D108").Value = "수면상태 무효화/수면상태의 시간이 2배가 된다."

$ws.Range("A109").Value = "수면병추가(睡眠瓶追加)"
$ws.Range("B109").Value = "수면병 추가(睡眠ビン追加)"
$ws.Range("D109").Value = "수면병의 장착이 가능하게 된다."

$ws.Range("A110").Value = "스태미너(スタミナ)"
$ws.Range("B110").Value = "러너(ランナー)/둔족(鈍足)"
$ws.Range("D110").Value = "대시, 귀인상태, 돌진, 차지 등의 스태미너 소비량이 1/2이 된다./대시, 귀인상태, 돌진, 차지 등의 스태미너 소비량이 1.2배가 된다."

$ws.Range("A111").Value = "청전주(青電主)"
$ws.Range("B111").Value = "청전주의 혼(青電主の魂)"
$ws.Range("D111").Value = "초회심(超会心) + 연발수+1(連発数+1) + 예리도 레벨+1(斬れ味レベル+1)"

$ws.Range("A112").Value = "정밀사격(精密射撃)"
$ws.Range("B112").Value = "흔들림 억제+2(ブレ抑制+2)/흔들림 억제+1(ブレ抑制+1)/흔들림 억제-1(ブレ抑制-1)/흔들림 억제-2(ブレ抑制-2)"
$ws.Range("D112").Value = "탄도의 흔들림 폭이 없어짐/탄도의 흔들림 폭이 1단계 감소한다./탄도의 흔들림 폭이 1단계 증가됨/탄도의 흔들림 폭이 대가 됨"

$ws.Range("A113").Value = "척안(隻眼)"
$ws.Range("B113").Value = "척안의 혼(隻眼の魂)"
$ws.Range("D113").Value = "기절무효(気絶無効) + 도전자+2(挑戦者+2)"

$ws.Range("A114").Value = "접격병추가(接撃瓶追加)"
$ws.Range("B114").Value = "접격병 추가(接撃ビン追加)"
$ws.Range("D114").Value = "접격병의 장착이 가능하게 된다."

$ws.Range("A115").Value = "절식(節食)"
$ws.Range("B115").Value = "만족감(満足感)"
$ws.Range("D115").Value = "먹거나 마시는 아이템이 25%의 확률로 소비하지 않게 된다."

$ws.Range("A116").Value = "천리안(千里眼)"
$ws.Range("B116").Value = "자동 마킹(自動マーキング)/탐지(探知)"
$ws.Range("D116").Value = "항상 중형, 대형 몬스터의 위치와 방향과 공격상태가 표시된다./페인트 시에 중형, 대형 몬스터의 위치와 방향과 공격상태가 표시된다. 또, 페인트 하지 않았을 때도 일정 확률로 같은 효과를 얻을 수 있다."

$ws.Range("A117").Value = "주행계속(走行継続)"
$ws.Range("B117").Value = "킵 런(キープラン)"
$ws.Range("D117").Value = "스태미너 잔량이 25(게이지 적색)이하에서 스태미너 소비와 속도저하가 없어지고 대시를 계속할 수 있게 된다."

$ws.Range("A118").Value = "장전수(装填数)"
$ws.Range("B118").Value = "장전수UP(装填数UP)"
$ws.Range("D118").Value = "보우건의 장전수나 활의 차지 단계가 하나 증가된다(최대 4단계까지). 건랜스의 포격이나 차지 액스의 장전수도 늘어난다."

$ws.Range("A119").Value = "장전속도(装填速度)"
$ws.Range("B119").Value = "장전속도+3(装填速度+3)/장전속도+2(装填速度+2)/장전속도+1(装填速度+1)/장전속도-1(装填速度-1)/장전속도-2(装填速度-2)/장전속도-3(装填速度-3)"
$ws.Range("D119").Value = "보우건의 리로드가 3단계 빨라짐. 앉아쏘기에는 영향없음. 활의 병이 교체하는 것만으로 자동으로 장전됨/보우건의 리로드가 2단계 빨라짐. 앉아쏘기에는 영향없음. 활의 병이 교체하는 것만으로 자동으로 장전됨./보우건의 리로드가 1단계 빨라짐. 앉아쏘기에는 영향 없음. 활의 병 교체시간이 0.75배로 빨라짐/보우건의 리로드가 1단계 느려짐. 앉아쏘기에는 영향없음. 활의 병 교체시간이 1.1배로 느려짐/보우건의 리로드가 2단계 느려짐. 앉아쏘기에는 영향없음. 활의 병 교체시간이 1.2배로 느려짐/보우건의 리로드가 3단계 느려짐. 앉아쏘기에는 영향없음. 활의 병 교체시간이 1.3배로 느려짐."

$ws.Range("A120").Value = "증폭(増幅)"
$ws.Range("B120").Value = "증폭강화(属物強化)"
$ws.Range("D120").Value = "강속성 공격(属性攻撃強化) + 아이템 사용 강화(アイテム使用強化)"

$ws.Range("A121").Value = "속강병추가(属強瓶追加)"
$ws.Range("B121").Value = "강속병 전LV 추가(属強ビン全LV追加)/강속병 LV1 추가(属強ビンLV1追加)"
$ws.Range("D121").Value = "전LV의 강속병을 사용할 수 있게 된다./강속병 LV1의 장착이 가능하게 된다."

$ws.Range("A122").Value = "속사(速射)"
$ws.Range("B122").Value = "연발수+1(連発数+1)"
$ws.Range("D122").Value = "속사시의 연발수가+1 된다."

$ws.Range("A123").Value = "속성회심(属性会心)"
$ws.Range("B123").Value = "회심격【속성】(会心撃【属性】)"
$ws.Range("D123").Value = "크리티컬 공격 시 가하는 속성 데미지(화, 수, 뇌, 빙, 용)를 높힌다. 대검은 1.2배.  보우건은 1.3배. 한손검, 쌍검, 활은 1.35배. 그 외는 1.25배"

$ws.Range("A124").Value = "속성공격(属性攻撃)"
$ws.Range("B124").Value = "속성공격강화(属性攻撃強化)/속성공격약화(属性攻撃弱化)"
$ws.Range("D124").Value = "화, 수, 뇌, 빙, 용의 속성치가 1.1배가 된다./화, 수, 뇌, 빙, 용의 속성치가 0.9가 된다."

$ws.Range("A125").Value = "속성내성(属性耐性)"
$ws.Range("B125").Value = "속성 피해 무효(属性やられ無効)"
$ws.Range("D125").Value = "화, 수, 뇌, 빙, 용의 속성 피해 무효화"

$ws.Range("A126").Value = "저력(底力)"
$ws.Range("B126").Value = "화사장력+2(火事場力+2)/화사장력+1(火事場力+1)/걱정이태산(心配性)"
$ws.Range("D126").Value = "체력이 40%이하가 되면 방어력이 45증가. 공격력 1.3배/체력이 40%이하가 되면 방어력이 45증가./체력이 40%이하가 되면 방어력이 +30에서 +21로 감소. 공격력 0.7배"

$ws.Range("A127").Value = "대염룡(対炎龍)"
$ws.Range("B127").Value = "강각의 수호(鋼殻の護り)"
$ws.Range("D127").Value = "남풍의 사냥꾼(南風の狩人) + 화내성【대】(火耐性【大】) + 세균 연구가(細菌研究家)"

$ws.Range("A128").Value = "대하룡(対霞龍)"
$ws.Range("B128").Value = "염린의 수호(炎鱗の護り)"
$ws.Range("D128").Value = "독내성(毒耐性) + 도난 무효(盗み無効) + 자동 마킹(自動マーキング)"

$ws.Range("A129").Value = "내한(耐寒)"
$ws.Range("B129").Value = "추위 무효(寒さ無効)/추위 배가(寒さ倍加)"
$ws.Range("D129").Value = "추위에 의한 스태미너 감소를 무효화한다./추위에 의한 스태미너 감소가 배가된다."

$ws.Range("A130").Value = "대강룡(対鋼龍)"
$ws.Range("B130").Value = "하피의 수호(霞皮の護り)"
$ws.Range("D130").Value = "북풍의 사냥꾼(北風の狩人) + 풍압【대】 무효(風圧【大】無効) + 눈사람 무효(だるま無効)"

$ws.Range("A131").Value = "체술(体術)"
$ws.Range("B131").Value = "체술+2(体術+2)/체술+1(体術+1)/체술-1(体術-1)/체술-2(体術-2)"
$ws.Range("D131").Value = "회피와 가드 등의 스태미너 소비가 절반이 된다. 대형 몬스터를 향해 긴급회피가 가능하게 되고 긴급회피의 이동거리가 늘어난다./회피와 가드 등의 스태미너 소비가 0.75배가 된다. 대형 몬스터를 향해 긴급회피가 가능하게 되고 긴급회피의 이동거리가 늘어난다./회피와 가드 등의 스태미너 소비가 1.2배로 증가한다./회피와 가드 등의 스태미너 소비가 1.35배로 증가한다."

$ws.Range("A132").Value = "내서(耐暑)"
$ws.Range("B132").Value = "더위 무효(暑さ無効)/더위 배가(暑さ倍加)"
$ws.Range("D132").Value = "더위, 용암지형, 불길에 의한 데미지를 무효화한다./더위, 용암지형, 불길에 의한 데미지가 배가된다."

$ws.Range("A133").Value = "내진(耐震)"
$ws.Range("B133").Value = "내진(耐震)"
$ws.Range("D133").Value = "진동에 의해 휘청거리는 상태를 무효화"

$ws.Range("A134").Value = "방어력DOWN(対防御DOWN)"
$ws.Range("B134").Value = "철면피(鉄面皮)"
$ws.Range("D134").Value = "방어력DOWN상태를 무효화"

$ws.Range("A135").Value = "체력(体力)"
$ws.Range("B135").Value = "체력+50(体力+50)/체력+20(体力+20)/체력-10(体力-10)/체력-30(体力-30)"
$ws.Range("D135").Value = "체력 최대치+50/체력 최대치+20/체력 최대치-10/체력 최대치-30"

$ws.Range("A136").Value = "보전(宝纏)"
$ws.Range("B136").Value = "보전의 혼(宝纏の魂)"
$ws.Range("D136").Value = "부적 헌터(お守りハンター) + 배고픔 무효(腹減り無効)"

$ws.Range("A137").Value = "장(匠)"
$ws.Range("B137").Value = "예리도 레벨+2(斬れ味レベル+2)/예리도 레벨+1(斬れ味レベル+1)"
$ws.Range("D137").Value = "예리도 게이지가 2단계 늘어난다./예리도 게이지가 1단계 늘어난다."

$ws.Range("A138").Value = "이식(茸食)"
$ws.Range("B138").Value = "버섯 애호가(キノコ大好き)"
$ws.Range("D138").Value = "버섯을 먹는 것으로 유효한 효과를 얻을 수 있다. 파란 버섯(アオキノコ)：회복약, 니트로 버섯(ニトロダケ)과 도스 송이버섯(ドスマツタケ)：귀인약, 마비버섯(マヒダケ)과 큰 마비 시메지(オオマヒシメジ)：경화약, 독광대 버섯(毒テングダケ)：영양제, 두근두근 버섯(ドキドキノコ)：랜덤이지만 나쁜 효과는 나오지 않음, 녹초버섯(クタビレダケ)：강주약, 엄선 버섯(厳選キノコ)：강주약 그레이트, 만드라고라(マンドラゴラ)：비약, 특산 버섯(特産キノコ)：휴대용 식량, 숙성 버섯(熟成キノコ)：고대의 비약, 심층 시메지(深層シメジ)：지움 열매, 걸상 버섯(コシカケダケ)：해독약, 혼돈 버섯(混沌茸)：천리안의 약"

$ws.Range("A139").Value = "달인(達人)"
$ws.Range("B139").Value = "간파+3(見切り+3)/간파+2(見切り+2)/간파+1(見切り+1)/간파-1(見切り-1)/간파-2(見切り-2)/간파-3(見切り-3)"
$ws.Range("D139").Value = "회심률+30%/회심률+20%/회심률+10%/회심률-5%/회심률-10%/회심률-15%"

$ws.Range("A140").Value = "순지(盾持)"
$ws.Range("B140").Value = "방패사용(盾使い)"
$ws.Range("D140").Value = "가드 강화(ガード強化) + 스태미너 급속회복(スタミナ急速回復)"

$ws.Range("D140").Select()
